# Auto-generated edit script applying the Jenova_Profits data refresh diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value2 = 380
$ws.Range("I18").Value2 = 380
$ws.Range("K18").Value2 = 380
$ws.Range("M18").Value2 = -96

$ws.Range("H19").Value2 = 366.5
$ws.Range("I19").Value2 = 118.625
$ws.Range("K19").Value2 = 118.625
$ws.Range("M19").Value2 = 56.375

$ws.Range("H125").Value2 = 11115066
$ws.Range("I125").Value2 = 848
$ws.Range("J125").Value2 = 15878302
$ws.Range("K125").Value2 = 7632
$ws.Range("L125").Value2 = 142904718
$ws.Range("M125").Value2 = -5172
$ws.Range("N125").Value2 = -142909638

$ws.Range("H132").Value2 = 2792.25
$ws.Range("I132").Value2 = 2426
$ws.Range("K132").Value2 = 7278
$ws.Range("M132").Value2 = -4748

$ws.Range("H137").Value2 = 3390.8096
$ws.Range("J137").Value2 = 4981.9287
$ws.Range("L137").Value2 = 14945.7861
$ws.Range("N137").Value2 = -20045.7861

$ws.Range("H138").Value2 = 6716.76
$ws.Range("I138").Value2 = 4026.8572
$ws.Range("J138").Value2 = 7334.1147
$ws.Range("K138").Value2 = 12080.5716
$ws.Range("L138").Value2 = 22002.3441
$ws.Range("M138").Value2 = -6940.571599999999
$ws.Range("N138").Value2 = -32282.3441

$ws.Range("H141").Value2 = 6208.2915
$ws.Range("I141").Value2 = 6186.3184
$ws.Range("K141").Value2 = 18558.9552
$ws.Range("M141").Value2 = -13378.9552

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value2 = 2651.5715
$ws.Range("I45").Value2 = 2001.2858
$ws.Range("K45").Value2 = 2001.2858
$ws.Range("M45").Value2 = -1624.2858

$ws.Range("H61").Value2 = 4658.385
$ws.Range("I61").Value2 = 4658.385
$ws.Range("K61").Value2 = 4658.385
$ws.Range("M61").Value2 = -4446.385

$ws.Range("H74").Value2 = 2101.4883
$ws.Range("I74").Value2 = 2114.842
$ws.Range("J74").Value2 = 2000
$ws.Range("K74").Value2 = 2114.842
$ws.Range("L74").Value2 = 2000
$ws.Range("M74").Value2 = -1240.842
$ws.Range("N74").Value2 = -3748

$ws.Range("H77").Value2 = 2101.4883
$ws.Range("I77").Value2 = 2114.842
$ws.Range("J77").Value2 = 2000
$ws.Range("K77").Value2 = 10574.21
$ws.Range("L77").Value2 = 10000
$ws.Range("M77").Value2 = -6206.210000000001
$ws.Range("N77").Value2 = -18736

$ws.Range("H122").Value2 = 3666.923
$ws.Range("I122").Value2 = 2476.818
$ws.Range("K122").Value2 = 7430.454000000001
$ws.Range("M122").Value2 = -4980.454000000001

$ws.Range("H136").Value2 = 4658.385
$ws.Range("I136").Value2 = 4658.385
$ws.Range("K136").Value2 = 13975.155
$ws.Range("M136").Value2 = -11425.155

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value2 = 49999.832
$ws.Range("J132").Value2 = 49999.832
$ws.Range("L132").Value2 = 49999.832
$ws.Range("N132").Value2 = -60119.832

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value2 = 728.3570999999999
$ws.Range("I22").Value2 = 426.9091
$ws.Range("J22").Value2 = 1833.6666
$ws.Range("K22").Value2 = 426.9091
$ws.Range("L22").Value2 = 1833.6666
$ws.Range("M22").Value2 = -76.90910000000002
$ws.Range("N22").Value2 = -2533.6666

$ws.Range("H31").Value2 = 30650.03
$ws.Range("I31").Value2 = 1545.619
$ws.Range("J31").Value2 = 74306.64
$ws.Range("K31").Value2 = 1545.619
$ws.Range("L31").Value2 = 74306.64
$ws.Range("M31").Value2 = -1250.619
$ws.Range("N31").Value2 = -74896.64

$ws.Range("H34").Value2 = 30650.03
$ws.Range("I34").Value2 = 1545.619
$ws.Range("J34").Value2 = 74306.64
$ws.Range("K34").Value2 = 1545.619
$ws.Range("L34").Value2 = 74306.64
$ws.Range("M34").Value2 = -1343.619
$ws.Range("N34").Value2 = -74710.64

$ws.Range("H58").Value2 = 2737.0952
$ws.Range("I58").Value2 = 2756.6667
$ws.Range("J58").Value2 = 2688.1667
$ws.Range("K58").Value2 = 2756.6667
$ws.Range("L58").Value2 = 2688.1667
$ws.Range("M58").Value2 = -2553.6667
$ws.Range("N58").Value2 = -3094.1667

$ws.Range("H99").Value2 = 6741.5454
$ws.Range("I99").Value2 = 6272.4
$ws.Range("K99").Value2 = 6272.4
$ws.Range("M99").Value2 = -4774.4

$ws.Range("H126").Value2 = 6741.5454
$ws.Range("I126").Value2 = 6272.4
$ws.Range("K126").Value2 = 18817.2
$ws.Range("M126").Value2 = -16347.2

$ws.Range("H132").Value2 = 3742.1304
$ws.Range("I132").Value2 = 2267.6875
$ws.Range("J132").Value2 = 7112.2856
$ws.Range("K132").Value2 = 6803.0625
$ws.Range("L132").Value2 = 21336.8568
$ws.Range("M132").Value2 = -4273.0625
$ws.Range("N132").Value2 = -26396.8568

$ws.Range("H136").Value2 = 2737.0952
$ws.Range("I136").Value2 = 2756.6667
$ws.Range("J136").Value2 = 2688.1667
$ws.Range("K136").Value2 = 8270.000100000001
$ws.Range("L136").Value2 = 8064.500100000001
$ws.Range("M136").Value2 = -5720.000100000001
$ws.Range("N136").Value2 = -13164.5001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value2 = 53809.316
$ws.Range("I129").Value2 = 614
$ws.Range("J129").Value2 = 84839.914
$ws.Range("K129").Value2 = 1842
$ws.Range("L129").Value2 = 254519.742
$ws.Range("M129").Value2 = 3158
$ws.Range("N129").Value2 = -264519.742

$ws.Range("H131").Value2 = 33156.523
$ws.Range("I131").Value2 = 84674
$ws.Range("J131").Value2 = 21492.19
$ws.Range("K131").Value2 = 254022
$ws.Range("L131").Value2 = 64476.56999999999
$ws.Range("M131").Value2 = -248982
$ws.Range("N131").Value2 = -74556.56999999999

$ws.Range("H134").Value2 = 2266.5186
$ws.Range("J134").Value2 = 9750
$ws.Range("L134").Value2 = 29250
$ws.Range("N134").Value2 = -39390

$ws.Range("H139").Value2 = 5828.1465
$ws.Range("I139").Value2 = 2916.5625
$ws.Range("J139").Value2 = 7691.56
$ws.Range("K139").Value2 = 8749.6875
$ws.Range("L139").Value2 = 23074.68
$ws.Range("M139").Value2 = -3609.6875
$ws.Range("N139").Value2 = -33354.68

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value2 = 4259.8
$ws.Range("I122").Value2 = 3791.6155
$ws.Range("K122").Value2 = 11374.8465
$ws.Range("M122").Value2 = -8924.8465

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 4878.357
$ws.Range("I7").Value2 = 3500.2856
$ws.Range("J7").Value2 = 6256.4287
$ws.Range("K7").Value2 = 3500.2856
$ws.Range("L7").Value2 = 6256.4287
$ws.Range("M7").Value2 = -3388.2856
$ws.Range("N7").Value2 = -6480.4287

$ws.Range("H40").Value2 = 116165.78
$ws.Range("I40").Value2 = 251448.5
$ws.Range("K40").Value2 = 251448.5
$ws.Range("M40").Value2 = -251312.5

$ws.Range("H126").Value2 = 4878.357
$ws.Range("I126").Value2 = 3500.2856
$ws.Range("J126").Value2 = 6256.4287
$ws.Range("K126").Value2 = 10500.8568
$ws.Range("L126").Value2 = 18769.2861
$ws.Range("M126").Value2 = -8030.856800000001
$ws.Range("N126").Value2 = -23709.2861

$ws.Range("H132").Value2 = 6660.864
$ws.Range("I132").Value2 = 5875.8887
$ws.Range("J132").Value2 = 7204.3076
$ws.Range("K132").Value2 = 17627.6661
$ws.Range("L132").Value2 = 21612.9228
$ws.Range("M132").Value2 = -15097.6661
$ws.Range("N132").Value2 = -26672.9228

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value2 = 600
$ws.Range("J13").Value2 = 0
$ws.Range("L13").Value2 = 0
$ws.Range("N13").ClearContents()

$ws.Range("H126").Value2 = 1832.6666
$ws.Range("I126").Value2 = 1776.8889
$ws.Range("J126").Value2 = 2000
$ws.Range("K126").Value2 = 5330.6667
$ws.Range("L126").Value2 = 6000
$ws.Range("M126").Value2 = -2860.6667
$ws.Range("N126").Value2 = -10940

$ws.Range("H132").Value2 = 1000000
$ws.Range("I132").Value2 = 0
$ws.Range("J132").Value2 = 1000000
$ws.Range("K132").Value2 = 0
$ws.Range("L132").Value2 = 3000000
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value2 = -3005060
